# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45186 (2023-09-17) to 45188 (2023-09-19).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
